# Componentenlijst RelaisPCB update
# "Update componentenlijst. Alle flyback diodes geintegreerd in de PCB desgin."
#  - flyback diode part swapped for a Diodes Incorporated 1N4448WQ-7-F, price drops to 0.2
#  - two new M3 hardware lines (screw + nut) sourced from Farnell added below it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: Fly back diode -> new part / new price ---
$ws.Range("B12").Value2 = 0.2

$newDiodeUrl = "https://www.mouser.be/ProductDetail/Diodes-Incorporated/1N4448WQ-7-F?qs=sGAEpiMZZMtoHjESLttvkiVPmB1TVWDoWOW8mzAF96J3zPSzsyGFzg%3D%3D"
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Row -eq 12) {
        $hl.Address = $newDiodeUrl
    }
}
$ws.Range("E12").Value2 = $newDiodeUrl

# --- Prep rows 13 & 14 with the same look (fonts/number formats) as row 12 ---
$ws.Range("A12:E12").Copy() | Out-Null
$ws.Range("A13:E13").PasteSpecial(-4122) | Out-Null
$ws.Range("A12:E12").Copy() | Out-Null
$ws.Range("A14:E14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- new component names first (keeps shared-string order: schroef, nut, urls, Farnell) ---
$ws.Range("A13").Value2 = "m3 schroef"
$ws.Range("A14").Value2 = "m3 nut"

$screwUrl = "https://be.farnell.com/tr-fastenings/m3-12-prstmc-z100/screw-pozi-pan-steel-bzp-m3x12/dp/1420391?ost=M3+12+PRSTMC+Z100&ddkey=https%3Anl-BE%2FElement14_Belgium%2Fsearch"
$ws.Range("E13").Value2 = $screwUrl
$ws.Hyperlinks.Add($ws.Range("E13"), $screwUrl) | Out-Null

$nutUrl = "https://be.farnell.com/tr-fastenings/m3-hfa2-s100/full-nut-stainless-steel-a2-m3/dp/1420788?ost=M3+-+HFA2-+S100&ddkey=https%3Anl-BE%2FElement14_Belgium%2Fsearch"
$ws.Range("E14").Value2 = $nutUrl
$ws.Hyperlinks.Add($ws.Range("E14"), $nutUrl) | Out-Null

# --- Row 13: m3 schroef (Farnell) ---
$ws.Range("B13").Value2 = 2.0299999999999998
$ws.Range("C13").Value2 = 1
$ws.Range("D13").Value2 = "Farnell"

# --- Row 14: m3 nut (Farnell) ---
$ws.Range("B14").Value2 = 3.48
$ws.Range("C14").Value2 = 1
$ws.Range("D14").Value2 = "Farnell"

# Re-apply formatting to the new hyperlink cells so they keep matching row 12's style
$ws.Range("E12").Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E13").Value2 = $screwUrl
$ws.Range("E14").Value2 = $nutUrl

# --- widen column E to fit the longer URLs ---
$ws.Columns.Item(5).ColumnWidth = 171.16666666666666

# --- view bookkeeping (zoom + selection) ---
$ws.Range("D22").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100
